{"js": "// Apply the documented edits to the \"Anwendung\" section of the document:\n//   1) \"sollten sie dies Wie im Beispiel\" -> \"sollten Sie dies wie im Beispiel\"\n//      (swap the capitalisation between \"sie\"/\"Sie\" and \"Wie\"/\"wie\")\n//   2) the curly quotes wrapping the comma change from single quotes\n//      (U+2018 ... U+2018) to double quotes (U+201D ... U+201C)\n//   3) the \"_GoBack\" bookmark moves from just before \"m\u00f6chten\" to just\n//      before \"weiterhin\" later in the document (after \"Fehler werden \")\n\nconst body = context.document.body;\n\n// --- Change 1 -----------------------------------------------------------\nlet results = body.search(\"sollten sie dies Wie im Beispiel\", { matchCase: true });\nresults.load(\"items\");\nawait context.sync();\n\nif (results.items.length > 0) {\n  results.items[0].insertText(\"sollten Sie dies wie im Beispiel\", \"Replace\");\n  await context.sync();\n}\n\n// --- Change 2 -------------------------------------------------------------\n// Search only for the punctuation run \"-\u2018,\u2018\" (not including the preceding\n// \"sep\" word) so the surrounding spell-check proofErr markers around \"sep\"\n// stay untouched.\nresults = body.search(\"-\\u2018,\\u2018\", { matchCase: true });\nresults.load(\"items\");\nawait context.sync();\n\nif (results.items.length > 0) {\n  results.items[0].insertText(\"-\\u201D,\\u201C\", \"Replace\");\n  await context.sync();\n}\n\n// --- Change 3: relocate the _GoBack bookmark -------------------------------\ncontext.document.deleteBookmark(\"_GoBack\");\nawait context.sync();\n\nresults = body.search(\"weiterhin als Error\", { matchCase: true });\nresults.load(\"items\");\nawait context.sync();\n\nif (results.items.length > 0) {\n  results.items[0].getRange(\"Start\").insertBookmark(\"_GoBack\");\n  await context.sync();\n}\n", "ps1": "# wdReplaceOne = 1 (replace the first/only match the Find locates)\n$wdReplaceOne = 1\n$wdFindContinue = 1\n\n$d = $word.ActiveDocument\n\n# --- Change 1: capitalisation swap -----------------------------------------\n# \"sollten sie dies Wie im Beispiel\" -> \"sollten Sie dies wie im Beispiel\"\n$rng1 = $d.Content\n$rng1.Find.ClearFormatting()\n$rng1.Find.Execute(\"sollten sie dies Wie im Beispiel\", $true, $false, $false, $false, $false, $true, $wdFindContinue, $false, \"sollten Sie dies wie im Beispiel\", $wdReplaceOne)\n\n# --- Change 2: quote style around the comma ---------------------------------\n# Search only the punctuation run \"-\u2018,\u2018\" (not the preceding \"sep\"\n# word) so the spell-check proofErr markers wrapping \"sep\" stay untouched.\n$rng2 = $d.Content\n$rng2.Find.ClearFormatting()\n$rng2.Find.Execute(\"-\u2018,\u2018\", $true, $false, $false, $false, $false, $true, $wdFindContinue, $false, \"-\u201d,\u201c\", $wdReplaceOne)\n\n# --- Change 3: relocate the _GoBack bookmark --------------------------------\nif ($d.Bookmarks.Exists(\"_GoBack\")) {\n    $d.Bookmarks(\"_GoBack\").Delete()\n}\n\n$rng3 = $d.Content\n$found3 = $rng3.Find.Execute(\"weiterhin als Error\", $true)\nif ($found3) {\n    $insertPoint = $d.Range($rng3.Start, $rng3.Start)\n    $d.Bookmarks.Add(\"_GoBack\", $insertPoint)\n}\n"}
